$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-22 Sunday" "2024-09-23 Monday"

Replace-Text "13×75=" "96×41="
Replace-Text "96×85=" "53×54="
Replace-Text "59×70=" "37×71="
Replace-Text "57×62=" "98×87="
Replace-Text "99×29=" "46×33="

Replace-Text "81×37=" "24×66="
Replace-Text "93×52=" "62×77="
Replace-Text "37×19=" "51×52="
Replace-Text "95×96=" "58×85="
Replace-Text "48×72=" "63×92="

Replace-Text "29×38=" "94×95="
Replace-Text "33×88=" "92×66="
Replace-Text "90×45=" "27×32="
Replace-Text "56×33=" "30×84="
Replace-Text "40×28=" "87×94="

Replace-Text "14×13=" "25×95="
Replace-Text "74×90=" "82×74="
Replace-Text "14×94=" "84×71="
Replace-Text "73×33=" "37×17="
Replace-Text "52×85=" "41×73="

Replace-Text "97×35=" "31×35="
Replace-Text "30×60=" "99×90="
Replace-Text "87×91=" "11×62="
Replace-Text "33×37=" "19×22="
Replace-Text "66×67=" "87×91="
